$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 10557.75
$ws.Range("I28").Value = 1043.1666
$ws.Range("J28").Value = 39101.5
$ws.Range("K28").Value = 1043.1666
$ws.Range("L28").Value = 39101.5
$ws.Range("M28").Value = -558.1666
$ws.Range("N28").Value = -40071.5

$ws.Range("H129").Value = 1528.65
$ws.Range("I129").Value = 583.6923
$ws.Range("K129").Value = 1751.0769
$ws.Range("M129").Value = 3248.9231

$ws.Range("H132").Value = 2807.1404
$ws.Range("I132").Value = 2712.16
$ws.Range("J132").Value = 3485.5715
$ws.Range("K132").Value = 8136.48
$ws.Range("L132").Value = 10456.7145
$ws.Range("M132").Value = -5606.48
$ws.Range("N132").Value = -15516.7145

$ws.Range("H137").Value = 1388.2069
$ws.Range("I137").Value = 1229.25
$ws.Range("J137").Value = 1648.3182
$ws.Range("K137").Value = 3687.75
$ws.Range("L137").Value = 4944.9546
$ws.Range("M137").Value = -1137.75
$ws.Range("N137").Value = -10044.9546

$ws.Range("H138").Value = 1860.87
$ws.Range("I138").Value = 1198.8125
$ws.Range("J138").Value = 1986.9762
$ws.Range("K138").Value = 3596.4375
$ws.Range("L138").Value = 5960.9286
$ws.Range("M138").Value = 1543.5625
$ws.Range("N138").Value = -16240.9286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1434.9048
$ws.Range("I2").Value = 1417.7778
$ws.Range("J2").Value = 1537.6666
$ws.Range("K2").Value = 1417.7778
$ws.Range("L2").Value = 1537.6666
$ws.Range("M2").Value = -1304.7778
$ws.Range("N2").Value = -1763.6666

$ws.Range("H32").Value = 792366.1
$ws.Range("I32").Value = 907282.75
$ws.Range("J32").Value = 19290.637
$ws.Range("K32").Value = 907282.75
$ws.Range("L32").Value = 19290.637
$ws.Range("M32").Value = -906995.75
$ws.Range("N32").Value = -19864.637

$ws.Range("H74").Value = 2073.16
$ws.Range("I74").Value = 1021.913
$ws.Range("J74").Value = 2968.6667
$ws.Range("K74").Value = 1021.913
$ws.Range("L74").Value = 2968.6667
$ws.Range("M74").Value = -147.913
$ws.Range("N74").Value = -4716.6667

$ws.Range("H77").Value = 2073.16
$ws.Range("I77").Value = 1021.913
$ws.Range("J77").Value = 2968.6667
$ws.Range("K77").Value = 5109.565000000001
$ws.Range("L77").Value = 14843.3335
$ws.Range("M77").Value = -741.5650000000005
$ws.Range("N77").Value = -23579.3335

$ws.Range("H110").Value = 53915.824
$ws.Range("I110").Value = 60972.6
$ws.Range("J110").Value = 990
$ws.Range("K110").Value = 60972.6
$ws.Range("L110").Value = 990
$ws.Range("M110").Value = -58927.6
$ws.Range("N110").Value = -5080

$ws.Range("H116").Value = 1434.9048
$ws.Range("I116").Value = 1417.7778
$ws.Range("J116").Value = 1537.6666
$ws.Range("K116").Value = 1417.7778
$ws.Range("L116").Value = 1537.6666
$ws.Range("M116").Value = 876.2221999999999
$ws.Range("N116").Value = -6125.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1434.9048
$ws.Range("I3").Value = 1417.7778
$ws.Range("J3").Value = 1537.6666
$ws.Range("K3").Value = 1417.7778
$ws.Range("L3").Value = 1537.6666
$ws.Range("M3").Value = -1303.7778
$ws.Range("N3").Value = -1765.6666

$ws.Range("H99").Value = 1050
$ws.Range("I99").Value = 910
$ws.Range("K99").Value = 910
$ws.Range("M99").Value = 588

$ws.Range("H133").Value = 44000
$ws.Range("J133").Value = 44000
$ws.Range("L133").Value = 44000
$ws.Range("N133").Value = -54120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4604.983
$ws.Range("I31").Value = 1425.3214
$ws.Range("K31").Value = 1425.3214
$ws.Range("M31").Value = -1130.3214

$ws.Range("H34").Value = 4604.983
$ws.Range("I34").Value = 1425.3214
$ws.Range("K34").Value = 1425.3214
$ws.Range("M34").Value = -1223.3214

$ws.Range("H94").Value = 1230.9474
$ws.Range("I94").Value = 850
$ws.Range("J94").Value = 1275.7646
$ws.Range("K94").Value = 850
$ws.Range("L94").Value = 1275.7646
$ws.Range("M94").Value = -399
$ws.Range("N94").Value = -2177.7646

$ws.Range("H132").Value = 2779071.8
$ws.Range("I132").Value = 1142.4445
$ws.Range("J132").Value = 11112860
$ws.Range("K132").Value = 3427.3335
$ws.Range("L132").Value = 33338580
$ws.Range("M132").Value = -897.3335000000002
$ws.Range("N132").Value = -33343640

$ws.Range("H134").Value = 3329.7114
$ws.Range("I134").Value = 3528.5898
$ws.Range("J134").Value = 2733.077
$ws.Range("K134").Value = 10585.7694
$ws.Range("L134").Value = 8199.231
$ws.Range("M134").Value = -8050.769400000001
$ws.Range("N134").Value = -13269.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1307.5
$ws.Range("I45").Value = 576.6667
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 1730.0001
$ws.Range("L45").Value = 10500
$ws.Range("M45").Value = -1198.0001
$ws.Range("N45").Value = -11564

$ws.Range("H115").Value = 5400.4287
$ws.Range("I115").Value = 3252.5
$ws.Range("J115").Value = 8264.333000000001
$ws.Range("K115").Value = 9757.5
$ws.Range("L115").Value = 24792.999
$ws.Range("M115").Value = -8582.5
$ws.Range("N115").Value = -27142.999

$ws.Range("H138").Value = 5975.933
$ws.Range("I138").Value = 1548.8889
$ws.Range("J138").Value = 12616.5
$ws.Range("K138").Value = 4646.6667
$ws.Range("L138").Value = 37849.5
$ws.Range("M138").Value = 493.3333000000002
$ws.Range("N138").Value = -48129.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 61500
$ws.Range("J68").Value = 98000
$ws.Range("L68").Value = 98000
$ws.Range("N68").Value = -99622

$ws.Range("H71").Value = 61500
$ws.Range("J71").Value = 98000
$ws.Range("L71").Value = 294000
$ws.Range("N71").Value = -302112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3996.2917
$ws.Range("I132").Value = 3631.5
$ws.Range("J132").Value = 4725.875
$ws.Range("K132").Value = 10894.5
$ws.Range("L132").Value = 14177.625
$ws.Range("M132").Value = -8364.5
$ws.Range("N132").Value = -19237.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9411207
$ws.Range("I132").Value = 3070.611
$ws.Range("J132").Value = 22437858
$ws.Range("K132").Value = 9211.832999999999
$ws.Range("L132").Value = 67313574
$ws.Range("M132").Value = -6681.832999999999
$ws.Range("N132").Value = -67318634
